# Update the numeric results in C2:H10 with corrected values
# (fixed minor mistake with construction of state variable; results changed only marginally)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(3.366303623731205, 12.15714401385543, 1.737188663037162, 4.004564280816123, 3.652049221412493, 24.91724980285242),
    @(4.686592407434943, 12.74007630712849, 1.067978942660879, 5.785149508426985, 4.629681100983577, 28.90947826663487),
    @(7.655722422417413, 11.33642055198433, 0.8619887241378588, 9.6321956492387, 4.461747400534877, 33.94807474831318),
    @(3.056804037521745, 11.28969770140409, 2.634483856901404, 3.383079511480429, 2.530826876317968, 22.89489198362563),
    @(4.419344227387025, 12.08477614690689, 1.72254167945185, 5.157664722330637, 3.462173094514133, 26.84649987059054),
    @(7.446834057065592, 10.79877460733823, 1.318006191099182, 9.236016582842852, 3.629350868742752, 32.42898230708861),
    @(3.510322822902935, 12.84694702643727, 3.168009596498625, 5.071380750934225, 4.196158676738227, 28.79281887351128),
    @(4.703387344963951, 13.3255840515091, 1.987108590042685, 6.903454256657646, 5.25270401146703, 32.17223825464041),
    @(7.33472612500473, 11.70590648267642, 1.709671881039287, 10.67236643672631, 5.045643801293241, 36.46831472673999)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowIndex = $i + 2
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colIndex = $j + 3  # Column C = 3
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
